$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column C ("Password") - this shifts D:J left to C:I,
# removing the Password / password1 / password2 / password3 values,
# and fixing the shared strings / styles / column widths accordingly.
$ws.Columns("C").Delete()

# Update the active selection to match the post-edit state.
$ws.Range("E9").Select()
